$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 (this shifts existing rows 33.. down by one,
# which is exactly the pattern shown by the diff: old row N (for N>=33) becomes
# row N+1 with identical content, and a brand-new record is inserted as the
# new row 33).
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record's data.
$ws.Cells.Item(33, 1).Value = 9
$ws.Cells.Item(33, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(33, 3).Value = "Metropolitana"
$ws.Cells.Item(33, 4).Value = 45028
$ws.Cells.Item(33, 5).Value = 13
$ws.Cells.Item(33, 6).Value = 100112029
$ws.Cells.Item(33, 7).Value = "Orégano"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 16
$ws.Cells.Item(33, 11).Value = 17000
$ws.Cells.Item(33, 12).Value = 17000
$ws.Cells.Item(33, 13).Value = 17000
$ws.Cells.Item(33, 14).Value = "$/docena de atados"
$ws.Cells.Item(33, 15).Value = "Región Metropolitana"
$ws.Cells.Item(33, 16).Value = 5667
$ws.Cells.Item(33, 17).Value = 3
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format style as the rest
# of column D.
$ws.Cells.Item(33, 4).NumberFormat = $ws.Cells.Item(34, 4).NumberFormat
